$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing metadata values in place -----------------------------
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Make room for the new "Jurisdiction" row (new row 12) ----------------
# Shift existing rows 12-15 (Description, Purpose, Copyright, Immutable)
# down to 13-16, then fill in row 12 with the new Jurisdiction property.
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Differential panel, method unspecified - Blood (69738-3)"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
